# Delete row 150 ("「お疲れ様です」..." post) entirely; all rows below shift up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(150).Delete()
